# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" and "全部类型" sheets, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 207
$ws1.Range("F4").Value  = 12838
$ws1.Range("F10").Value = 215
$ws1.Range("F11").Value = 464
$ws1.Range("F12").Value = 64
$ws1.Range("F16").Value = 394
$ws1.Range("F17").Value = 5474
$ws1.Range("F18").Value = 102
$ws1.Range("F19").Value = 38
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 129
$ws1.Range("F23").Value = 115

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 207
$ws4.Range("F4").Value  = 12839
$ws4.Range("F10").Value = 215
$ws4.Range("F11").Value = 464
$ws4.Range("F12").Value = 64
$ws4.Range("F16").Value = 394
$ws4.Range("F17").Value = 5474
$ws4.Range("F18").Value = 102
$ws4.Range("F19").Value = 38
$ws4.Range("F21").Value = 28
$ws4.Range("F22").Value = 129
$ws4.Range("F23").Value = 115

$wb.Save()
